# Add a new "buttonLabel" / "Record Location" column (F) to the survey
# sheet, driving the launch-intent button's text, and make "survey" the
# active/selected sheet (it was "settings" before).
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("survey")
$ws2 = $wb.Worksheets.Item("settings")

$ws1.Range("F1").Value = "buttonLabel"
$ws1.Range("F5").Value = "Record Location"

# New column F is a bit wider than the rest (~14.5 "characters" of stored
# OOXML width, which is ~0.833 wider than the ColumnWidth figure Excel
# reports back through COM).
$ws1.Columns.Item(6).ColumnWidth = 13.666666666666666

# Restore settings' prior selection before switching away from it, then
# make survey the active sheet/selection, matching the new view state.
$ws2.Range("C8").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("F6").Select() | Out-Null
